$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 16 for the new IPO record ("하스"), shifting rows 16-24 down to 17-25
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16
$ws.Cells.Item(16, 1).Value = "삼성"
$ws.Cells.Item(16, 2).Value = "2024-06-24"
$ws.Cells.Item(16, 3).Value = "하스"
$ws.Cells.Item(16, 4).Value = "삼성"
$ws.Cells.Item(16, 5).Value = "삼성"
$ws.Cells.Item(16, 6).Value = "2024-06-27"
$ws.Cells.Item(16, 7).Value = "2024-07-03"
$ws.Cells.Item(16, 8).Value = 28960
$ws.Cells.Item(16, 9).Value = 1810000
$ws.Cells.Item(16, 10).Value = 16000
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 100

# Apply data corrections reconciled in this RPA push
# row 7
$ws.Cells.Item(7, 2).Value = "2024-06-19"
$ws.Cells.Item(7, 3).Value = "에이치브이엠"
$ws.Cells.Item(7, 6).Value = "2024-06-24"
$ws.Cells.Item(7, 7).Value = "2024-06-28"
$ws.Cells.Item(7, 8).Value = 43200
$ws.Cells.Item(7, 9).Value = 2400000
$ws.Cells.Item(7, 10).Value = 18000
# row 8
$ws.Cells.Item(8, 2).Value = "2024-05-07"
$ws.Cells.Item(8, 3).Value = "아이씨티케이"
$ws.Cells.Item(8, 6).Value = "2024-05-10"
$ws.Cells.Item(8, 7).Value = "2024-05-17"
$ws.Cells.Item(8, 8).Value = 39400
$ws.Cells.Item(8, 9).Value = 1970000
$ws.Cells.Item(8, 10).Value = 20000
# row 10
$ws.Cells.Item(10, 2).Value = "2024-06-10"
$ws.Cells.Item(10, 3).Value = "미래에셋비전스팩5호"
$ws.Cells.Item(10, 5).Value = "미래"
$ws.Cells.Item(10, 6).Value = "2024-06-13"
$ws.Cells.Item(10, 7).Value = "2024-06-19"
$ws.Cells.Item(10, 8).Value = 9500
$ws.Cells.Item(10, 9).Value = 4750000
$ws.Cells.Item(10, 10).Value = 2000
$ws.Cells.Item(10, 12).Value = 100
# row 12
$ws.Cells.Item(12, 2).Value = "2024-06-13"
$ws.Cells.Item(12, 3).Value = "미래에셋비전스팩6호"
$ws.Cells.Item(12, 6).Value = "2024-06-18"
$ws.Cells.Item(12, 7).Value = "2024-06-24"
$ws.Cells.Item(12, 8).Value = 12900
$ws.Cells.Item(12, 9).Value = 6450000
# row 13
$ws.Cells.Item(13, 2).Value = "2024-06-20"
$ws.Cells.Item(13, 3).Value = "이노스페이스"
$ws.Cells.Item(13, 5).Value = "미래, 신한"
$ws.Cells.Item(13, 6).Value = "2024-06-25"
$ws.Cells.Item(13, 7).Value = "2024-07-02"
$ws.Cells.Item(13, 8).Value = 54133.66
$ws.Cells.Item(13, 9).Value = 1330000
$ws.Cells.Item(13, 10).Value = 43300
$ws.Cells.Item(13, 12).Value = 94
# row 19
$ws.Cells.Item(19, 3).Value = "한국제14호스팩"
$ws.Cells.Item(19, 8).Value = 8000
$ws.Cells.Item(19, 9).Value = 4000000
$ws.Cells.Item(19, 10).Value = 2000
# row 23
$ws.Cells.Item(23, 3).Value = "씨어스테크놀로지"
$ws.Cells.Item(23, 8).Value = 22100
$ws.Cells.Item(23, 9).Value = 1300000
$ws.Cells.Item(23, 10).Value = 17000
